$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44574
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101004
$ws.Cells.Item($row, 10).Value = "Frambuesa"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 200
$ws.Cells.Item($row, 14).Value = 3000
$ws.Cells.Item($row, 15).Value = 3000
$ws.Cells.Item($row, 16).Value = 3000
$ws.Cells.Item($row, 17).Value = "`$/envase 1 kilo"
$ws.Cells.Item($row, 18).Value = "Región de La Araucanía"
$ws.Cells.Item($row, 19).Value = 3000
$ws.Cells.Item($row, 20).Value = 1
